$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.472.68'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '3.044.95'
$ws.Range('E3').Value = '  +4.51%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '201.85'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.32'
$ws.Range('E6').Value = '  +4.72%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.551'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.208'
$ws.Range('E9').Value = '  +6.16%  '
$ws.Range('D10').Value = '3.044.98'
$ws.Range('E10').Value = '  +4.48%  '
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.11'
$ws.Range('E13').Value = '  +4.72%  '
$ws.Range('D14').Value = '3.608.71'
$ws.Range('E14').Value = '  +4.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.62'
$ws.Range('E15').Value = '  +6.67%  '
$ws.Range('D16').Value = '76.444.45'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('E17').Value = '  +1.94%  '
$ws.Range('D18').Value = '3.053.72'
$ws.Range('E18').Value = '  +4.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.49'
$ws.Range('E19').Value = '  +4.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.00'
$ws.Range('E20').Value = '  +1.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '374.65'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.30'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.35'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.51'
$ws.Range('E24').Value = '  +3.26%  '
$ws.Range('D25').Value = '3.204.40'
$ws.Range('E25').Value = '  +4.93%  '
$ws.Range('E26').Value = '  +4.11%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('E29').Value = '  +3.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +7.97%  '
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '512.93'
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('E34').Value = '  +7.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '20.89'
$ws.Range('E36').Value = '  +3.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '163.37'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.386'
$ws.Range('E38').Value = '  +8.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '20.03'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('E40').Value = '  +1.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '188.80'
$ws.Range('E41').Value = '  +4.28%  '
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.20'
$ws.Range('E44').Value = '  +4.54%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.789'
$ws.Range('E45').Value = '  +20.57%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.28'
$ws.Range('E46').Value = '  +7.53%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '42.20'
$ws.Range('E47').Value = '  +5.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.66'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.45'
$ws.Range('E49').Value = '  +4.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.607'
$ws.Range('E50').Value = '  +6.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.88'
$ws.Range('E51').Value = '  +4.89%  '
